# "Termino da atividade 'Prova 2des'" - recompute Lucro % as Lucro R$ / Faturamento
# (previously it duplicated the Lucro R$ formula by mistake), tidy up the
# leftover formatting below the table and restore the print/page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Lucro %" calculated column -----------------------------------
# It used to repeat Faturamento-Investimento (same as "Lucro R$"); it should
# really be the profit ratio: Lucro R$ / Faturamento.
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("I$r").Formula = "=Tabela1[[#This Row],[Lucro R$]]/Tabela1[[#This Row],[Faturamento]]"
}

# --- Clean up the empty formatted cells left under the table ---------------
$ws.Range("F12:H12").Clear()
$ws.Range("B12").Font.Underline = $true

# --- Restore the page setup used for printing -------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Move the selection/view to where the user left off ---------------------
$ws.Range("B12").Select() | Out-Null
